# Apply regression-results update: refreshed "Best1" stats (rows 3-13)
# and a newly added "Best2" model block (rows 15-26), per filter.py
# + correction function changes (see commit message).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Best1 block: updated summary stats + coefficient table ---
$ws.Range("E3").Value = "'0.729"

$ws.Range("E4").Value = "'29.9829"

$ws.Range("C5").Value = "2025-02-25 01:14"
$ws.Range("E5").Value = "'34.8584"

$ws.Range("E6").Value = "'-10.991"

$ws.Range("C7").Value = "'3"
$ws.Range("E7").Value = "'22.52"

$ws.Range("C8").Value = "'21"
$ws.Range("E8").Value = "'9.16e-07"

$ws.Range("C9").Value = "'0.763"
$ws.Range("E9").Value = "'0.16793"

$ws.Range("G10").Value = -4.488367156397095
$ws.Range("H10").Value = 0.9755853017823947
$ws.Range("I10").Value = -4.600691654739823
$ws.Range("J10").Value = 0.0001546352537605774
$ws.Range("K10").Value = -6.517207856696577
$ws.Range("L10").Value = -2.459526456097612

$ws.Range("F11").Value = "平均気温(℃)_4_下旬"
$ws.Range("G11").Value = 0.1915057525007337
$ws.Range("H11").Value = 0.05381629016720139
$ws.Range("I11").Value = 3.558508992458343
$ws.Range("J11").Value = 0.001857197708189514
$ws.Range("K11").Value = 0.07958865039714048
$ws.Range("L11").Value = 0.3034228546043268
$ws.Range("M11").Value = 0.3853869628902444

$ws.Range("F12").Value = "日最高気温の平均(℃)_1_中旬"
$ws.Range("G12").Value = 0.1656552634562546
$ws.Range("H12").Value = 0.05772178324513121
$ws.Range("I12").Value = 2.869891644767012
$ws.Range("J12").Value = 0.009168194065138335
$ws.Range("K12").Value = 0.04561624387731052
$ws.Range("L12").Value = 0.2856942830351987
$ws.Range("M12").Value = 0.312417841616266

$ws.Range("F13").Value = "日最低気温の平均(℃)_1_下旬"
$ws.Range("G13").Value = 0.43138041026291
$ws.Range("H13").Value = 0.06694508164267994
$ws.Range("I13").Value = 6.443795416747825
$ws.Range("J13").Value = [double]"2.189505814400155e-06"
$ws.Range("K13").Value = 0.2921604916423691
$ws.Range("L13").Value = 0.5706003288834508
$ws.Range("M13").Value = 0.7084335814712248

# --- blank spacer row 14 (kept as empty cells, like row 2) ---
$ws.Range("A14:M14").Style = "Normal"

# --- Best2 block: newly added second model ---
$ws.Range("A15").Value = "Best2"

$ws.Range("B16").Value = "Model:"
$ws.Range("C16").Value = "OLS"
$ws.Range("D16").Value = "Adj. R-squared:"
$ws.Range("E16").Value = "'0.687"

$ws.Range("B17").Value = "Dependent Variable:"
$ws.Range("C17").Value = "発病率"
$ws.Range("D17").Value = "AIC:"
$ws.Range("E17").Value = "'33.5517"

$ws.Range("B18").Value = "Date:"
$ws.Range("C18").Value = "2025-02-25 01:14"
$ws.Range("D18").Value = "BIC:"
$ws.Range("E18").Value = "'38.4272"

$ws.Range("B19").Value = "No. Observations:"
$ws.Range("C19").Value = "'25"
$ws.Range("D19").Value = "Log-Likelihood:"
$ws.Range("E19").Value = "'-12.776"

$ws.Range("B20").Value = "Df Model:"
$ws.Range("C20").Value = "'3"
$ws.Range("D20").Value = "F-statistic:"
$ws.Range("E20").Value = "'18.60"

$ws.Range("B21").Value = "Df Residuals:"
$ws.Range("C21").Value = "'21"
$ws.Range("D21").Value = "Prob (F-statistic):"
$ws.Range("E21").Value = "'4.01e-06"

$ws.Range("B22").Value = "R-squared:"
$ws.Range("C22").Value = "'0.727"
$ws.Range("D22").Value = "Scale:"
$ws.Range("E22").Value = "'0.19370"

$ws.Range("F23").Value = "const"
$ws.Range("G23").Value = -3.450274858657991
$ws.Range("H23").Value = 0.9379734148112036
$ws.Range("I23").Value = -3.678435661582654
$ws.Range("J23").Value = 0.001397879783702184
$ws.Range("K23").Value = -5.400897358085853
$ws.Range("L23").Value = -1.49965235923013

$ws.Range("F24").Value = "平均気温(℃)_1_中旬"
$ws.Range("G24").Value = 0.1764849364083666
$ws.Range("H24").Value = 0.08464786828717882
$ws.Range("I24").Value = 2.084930666057871
$ws.Range("J24").Value = 0.04946820310757518
$ws.Range("K24").Value = 0.0004500575916659144
$ws.Range("L24").Value = 0.3525198152250674
$ws.Range("M24").Value = 0.2646267102789007

$ws.Range("F25").Value = "平均気温(℃)_4_下旬"
$ws.Range("G25").Value = 0.1648783656538122
$ws.Range("H25").Value = 0.06115255368082748
$ws.Range("I25").Value = 2.696181201432063
$ws.Range("J25").Value = 0.01352387920233384
$ws.Range("K25").Value = 0.03770466837871178
$ws.Range("L25").Value = 0.2920520629289126
$ws.Range("M25").Value = 0.3318019002347538

$ws.Range("F26").Value = "日最低気温の平均(℃)_1_下旬"
$ws.Range("G26").Value = 0.4153406184311692
$ws.Range("H26").Value = 0.07520224956236303
$ws.Range("I26").Value = 5.522981305057096
$ws.Range("J26").Value = [double]"1.762889584589413e-05"
$ws.Range("K26").Value = 0.2589489790866143
$ws.Range("L26").Value = 0.5717322577757241
$ws.Range("M26").Value = 0.6820922666987538

# --- trailing blank spacer row 27 (kept as empty cells) ---
$ws.Range("A27:M27").Style = "Normal"
